# Woche 9 (KW51: 16.12. - 22.12.2024)
#
# Kosten.xlsx / Tabelle1 - update the cost-tracking table:
#   - remove the "Deans-T-Adapter" line item (row 11) - no longer purchased
#   - add two new flat-rate ("pauschal") line items at the bottom of the
#     item list: "3D-Druck Filament pauschal" and
#     "sonstige Verbrauchsmaterialien pauschal"
#   - totals / helper formulas re-flow automatically because they are
#     written relative to the (now one-row-shorter) data block

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# 1) Drop the "Deans-T-Adapter" row entirely; Excel shifts rows 12-20 up by
#    one, which also slides the two trailing blank rows and the totals row
#    into their new (A15:F19 range -> totals now on row 19) positions and
#    keeps every per-row style (borders/number formats) intact because they
#    are identical for every line-item row.
$ws.Rows.Item(11).Delete() | Out-Null

# 2) Fill the first now-available blank line (row 15) with the new 3D-print
#    filament flat-rate entry.
$ws.Range("A15").Value2 = "3D-Druck Filament pauschal"
$ws.Range("B15").Value2 = 1
$ws.Range("C15").Value2 = 20
$ws.Range("D15").Value2 = 1
$ws.Range("F15").Value2 = "Marcel List"
$ws.Range("E15").Formula = '=IF(B15>0,(C15/B15)*D15,"")'

# 3) Fill the second now-available blank line (row 16) with the
#    miscellaneous-consumables flat-rate entry, attributed to "alle" (all
#    team members) rather than a single person.
$ws.Range("A16").Value2 = "sonstige Verbrauchsmaterialien pauschal"
$ws.Range("B16").Value2 = 1
$ws.Range("C16").Value2 = 5
$ws.Range("D16").Value2 = 1
$ws.Range("F16").Value2 = "alle"
$ws.Range("E16").Formula = '=IF(B16>0,(C16/B16)*D16,"")'

# 3b) Re-assert the unaffected shared-formula cells too, so every cost cell
#     in the (now 18-row) data block is freshly (re)calculated.
$ws.Range("E2").Formula  = '=IF(B2>0,(C2/B2)*D2,"")'
$ws.Range("E17").Formula = '=IF(B17>0,(C17/B17)*D17,"")'
$ws.Range("E18").Formula = '=IF(B18>0,(C18/B18)*D18,"")'

# 4) Force a full recalculation so the SUM totals (now SUM(C2:C18) /
#    SUM(E2:E18)) and every per-row formula pick up the new data.
$excel.CalculateFullRebuild()

# 5) Keep the on-open selection sane (matches the saved worksheet state).
$ws.Range("M14").Select() | Out-Null
